$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number + week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/19/2024  Through  8/25/2024"

# --- Cells whose type/style flips between text ("0"/"***.*") and numeric ---
# Use same-row template cells (C14="0" text, E14="***.*" text, F14=count number,
# K14=percent number) as style donors via Copy, then overwrite the value where needed.
$ws.Range("C14").Copy($ws.Range("G14"))
$ws.Range("E14").Copy($ws.Range("H14"))

$ws.Range("F14").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("K14").Copy($ws.Range("E15"))
$ws.Range("E15").Value = 0
$ws.Range("F14").Copy($ws.Range("G15"))
$ws.Range("G15").Value = 1
$ws.Range("K14").Copy($ws.Range("H15"))
$ws.Range("H15").Value = 100

$ws.Range("F14").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 1

$ws.Range("C14").Copy($ws.Range("C20"))

$ws.Range("F14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("K14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100

$ws.Range("F14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("K14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = 100

$ws.Range("F14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1

$ws.Range("F14").Copy($ws.Range("C29"))
$ws.Range("C29").Value = 1
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))

$ws.Range("F14").Copy($ws.Range("C30"))
$ws.Range("C30").Value = 1
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))

$ws.Range("C14").Copy($ws.Range("C33"))

# --- Plain numeric value updates (style/format unchanged) ---
$ws.Range("N14").Value = -70.833333333333
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 8
$ws.Range("K15").Value = 12.5
$ws.Range("L15").Value = -10
$ws.Range("M15").Value = -25
$ws.Range("N15").Value = -65.384615384615
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 25
$ws.Range("G16").Value = 29
$ws.Range("H16").Value = -13.793103448275
$ws.Range("I16").Value = 174
$ws.Range("J16").Value = 171
$ws.Range("K16").Value = 1.754385964912
$ws.Range("L16").Value = 10.828025477707
$ws.Range("M16").Value = 25.179856115107
$ws.Range("N16").Value = -68.535262206148
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 233.333333333333
$ws.Range("F17").Value = 37
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 37.037037037037
$ws.Range("I17").Value = 327
$ws.Range("J17").Value = 263
$ws.Range("K17").Value = 24.334600760456
$ws.Range("L17").Value = 5.144694533762
$ws.Range("M17").Value = 108.28025477707
$ws.Range("N17").Value = -26.351351351351
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -28.571428571428
$ws.Range("I18").Value = 64
$ws.Range("K18").Value = -17.948717948717
$ws.Range("L18").Value = -54.609929078014
$ws.Range("M18").Value = -9.859154929577
$ws.Range("N18").Value = -80.780780780780
$ws.Range("C19").Value = 10
$ws.Range("E19").Value = 11.111111111111
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -10.638297872340
$ws.Range("I19").Value = 340
$ws.Range("J19").Value = 314
$ws.Range("K19").Value = 8.280254777070
$ws.Range("L19").Value = -1.449275362318
$ws.Range("M19").Value = 142.857142857143
$ws.Range("N19").Value = 28.787878787878
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 10
$ws.Range("H20").Value = 25
$ws.Range("J20").Value = 54
$ws.Range("K20").Value = 44.444444444444
$ws.Range("L20").Value = 44.444444444444
$ws.Range("M20").Value = 65.957446808510
$ws.Range("N20").Value = -63.720930232558
$ws.Range("C21").Value = 28
$ws.Range("E21").Value = 12
$ws.Range("F21").Value = 122
$ws.Range("G21").Value = 119
$ws.Range("H21").Value = 2.521008403361
$ws.Range("I21").Value = 999
$ws.Range("J21").Value = 894
$ws.Range("K21").Value = 11.744966442953
$ws.Range("L21").Value = -2.346041055718
$ws.Range("M21").Value = 74.041811846689
$ws.Range("N21").Value = -46.261430876815
$ws.Range("G22").Value = 2
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = -18.75
$ws.Range("L22").Value = -55.172413793103
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 14
$ws.Range("G23").Value = 11
$ws.Range("H23").Value = 27.272727272727
$ws.Range("I23").Value = 113
$ws.Range("J23").Value = 110
$ws.Range("K23").Value = 2.727272727272
$ws.Range("L23").Value = -3.418803418803
$ws.Range("M23").Value = 66.176470588235
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 64
$ws.Range("E24").Value = -71.875
$ws.Range("G24").Value = 128
$ws.Range("H24").Value = -27.34375
$ws.Range("I24").Value = 666
$ws.Range("J24").Value = 834
$ws.Range("K24").Value = -20.143884892086
$ws.Range("L24").Value = -21.462264150943
$ws.Range("M24").Value = 9.900990099009
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 49
$ws.Range("E25").Value = -89.795918367346
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 83
$ws.Range("H25").Value = -84.337349397590
$ws.Range("I25").Value = 135
$ws.Range("J25").Value = 403
$ws.Range("K25").Value = -66.501240694789
$ws.Range("L25").Value = -72.222222222222
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 350
$ws.Range("F26").Value = 63
$ws.Range("G26").Value = 27
$ws.Range("H26").Value = 133.333333333333
$ws.Range("I26").Value = 502
$ws.Range("J26").Value = 330
$ws.Range("K26").Value = 52.121212121212
$ws.Range("L26").Value = 38.292011019283
$ws.Range("M26").Value = 49.850746268656
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 14
$ws.Range("J27").Value = 14
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -30
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 5
$ws.Range("I28").Value = 41
$ws.Range("J28").Value = 36
$ws.Range("K28").Value = 13.888888888888
$ws.Range("L28").Value = -14.583333333333
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 12
$ws.Range("K29").Value = -47.826086956521
$ws.Range("L29").Value = -62.5
$ws.Range("M29").Value = -57.142857142857
$ws.Range("N29").Value = -74.468085106383
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -50
$ws.Range("I30").Value = 9
$ws.Range("K30").Value = -59.090909090909
$ws.Range("L30").Value = -59.090909090909
$ws.Range("M30").Value = -60.869565217391
$ws.Range("N30").Value = -80.434782608695
